# Generate Report for Handoff
# Rotate the handoff package identifiers: the source doc guid
# "e1110457-54ff-4265-a6b3-d7a950a6eb6f" becomes
# "66098900-5720-41b3-9af7-5ede14027edf", the xlf package hash
# "15eb4977bcaeec8a41783a258ceb618bd39cbddb" becomes
# "01d2803b7a81760899f7f0ec37de7c1927b748e3", and the handoff
# timestamps advance to the new run.

$wb = $excel.ActiveWorkbook

$oldGuid = "e1110457-54ff-4265-a6b3-d7a950a6eb6f"
$newGuid = "66098900-5720-41b3-9af7-5ede14027edf"
$oldHash = "15eb4977bcaeec8a41783a258ceb618bd39cbddb"
$newHash = "01d2803b7a81760899f7f0ec37de7c1927b748e3"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"

$oldZhXlf = "$oldGuid.$oldHash.zh-cn.xlf"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$oldDeXlf = "$oldGuid.$oldHash.de-de.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

$oldZhDate = "2016-03-04 06:20:52"
$newZhDate = "2016-03-04 06:22:01"
$oldDeDate = "2016-03-04 06:21:08"
$newDeDate = "2016-03-04 06:22:17"

# --- Overview sheet: just the source-doc filename + its hyperlink ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    }
}

# --- zh-cn sheet: source-doc filename, handoff xlf filename + its datetime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("C2").Value = $newZhXlf
$wsZh.Range("D2").Value = $newZhDate
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.TextToDisplay -eq $oldZhXlf) {
        $hl.TextToDisplay = $newZhXlf
    }
}

# --- de-de sheet: source-doc filename, handoff xlf filename + its datetime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("C2").Value = $newDeXlf
$wsDe.Range("D2").Value = $newDeDate
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.TextToDisplay -eq $oldDeXlf) {
        $hl.TextToDisplay = $newDeXlf
    }
}
